# The slide master and every slide layout carry a cached "last known value"
# for the auto-updating datetimeFigureOut field (the one shown in the
# Insert > Header & Footer "Date and time" box). The deck was re-saved on
# 2022-08-08, so that cached value needs to move from 2022-06-14 to
# 2022-08-08 everywhere it appears: once on the slide master and once on
# each of the eleven slide layouts.

$p = $ppt.ActivePresentation
$newDate = "2022-08-08"

function Set-DatePlaceholderText {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }

        $phType = -1
        try { $phType = $shp.PlaceholderFormat.Type } catch {}

        # 16 == ppPlaceholderDate (the "dt" placeholder holding the
        # datetimeFigureOut field).
        if ($phType -eq 16) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master.
$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes

# Every slide layout hanging off the master.
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Set-DatePlaceholderText $layout.Shapes
}
